# Staging.PeopleReachedValues.xlsx update:
# "staging templates updated with database changes - framework and impact"
#
# - Column header previously labelled "Framework_IndicatorBusinessKey" (G2)
#   is renamed to "FrameworkBusinessKey".
# - Column header previously labelled "FrameworkBusinessKey" (H2) is renamed
#   to "FrameworkDetail_IndicatorBusinessKey" (keeps the alphabetical sort
#   order: FrameworkBusinessKey < FrameworkDetail_IndicatorBusinessKey < Gender...).
# - A new trailing column "NumberReached" is appended after "Notes".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two Framework-related headers.
$ws.Range("G2").Value = "FrameworkBusinessKey"
$ws.Range("H2").Value = "FrameworkDetail_IndicatorBusinessKey"

# Add the new "NumberReached" header in the next empty column (W), copying
# the bold/underline header formatting used by the rest of row 2.
$ws.Range("V2").Copy($ws.Range("W2"))
$ws.Range("W2").Value = "NumberReached"
